$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.483.52"
$ws.Range("E2").Value = "  -0.40%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.716.28"
$ws.Range("E3").Value = "  -0.50%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.59"
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.31"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.715.39"
$ws.Range("E7").Value = "  -0.40%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  -2.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  -2.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.53"
$ws.Range("E11").Value = "  +2.10%  "

$ws.Range("E12").Value = "  -4.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.57"
$ws.Range("E13").Value = "  -3.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000253"
$ws.Range("E14").Value = "  -1.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.337.51"
$ws.Range("E15").Value = "  -0.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.719.38"
$ws.Range("E16").Value = "  -0.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.548.06"
$ws.Range("E17").Value = "  -0.32%  "

$ws.Range("E18").Value = "  -2.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.47"
$ws.Range("E19").Value = "  -1.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "500.07"
$ws.Range("E20").Value = "  -3.42%  "

$ws.Range("E21").Value = "  -3.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.13"
$ws.Range("E22").Value = "  -2.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.717"
$ws.Range("E23").Value = "  -2.09%  "

$ws.Range("E24").Value = "  +3.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.95"
$ws.Range("E25").Value = "  -2.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.20"
$ws.Range("E26").Value = "  +1.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.87"
$ws.Range("E27").Value = "  -5.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000134"
$ws.Range("E28").Value = "  +5.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("E30").Value = "  -3.01%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.90"
$ws.Range("E31").Value = "  +1.63%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.97"
$ws.Range("E32").Value = "  +1.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.29"
$ws.Range("E33").Value = "  -3.83%  "

$ws.Range("E34").Value = "  -2.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.04"
$ws.Range("E36").Value = "  +0.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.06"
$ws.Range("E37").Value = "  -2.68%  "

$ws.Range("E38").Value = "  +1.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.137"
$ws.Range("E39").Value = "  +4.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.01"
$ws.Range("E40").Value = "  +9.49%  "

$ws.Range("E41").Value = "  -6.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "49.65"
$ws.Range("E42").Value = "  -3.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "45.37"
$ws.Range("E43").Value = "  +1.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "434.21"
$ws.Range("E44").Value = "  +2.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.54"
$ws.Range("E45").Value = "  -3.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.946.28"
$ws.Range("E46").Value = "  -3.75%  "

$ws.Range("E47").Value = "  -1.40%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "139.25"
$ws.Range("E48").Value = "  +2.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.96"
$ws.Range("E50").Value = "  -3.74%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.46"
$ws.Range("E51").Value = "  -3.20%  "
